# PrefArr.xlsx edit: a new data row (45) was appended under the
# "Attributes" column on Sheet1, and the sheet's on-screen selection
# moved down to the first empty cell below the new data (B47), with the
# view scrolled so row 28 is at the top.
#
# Reproduce by writing the new value (extends the used range /
# <dimension> from A1:A44 to A1:A45 automatically) and then moving the
# active selection the same way a user would after typing the value and
# pressing Enter a few times.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data point in column A, row 45 (matches the existing
# plain (unstyled) cells A36:A44 - no explicit style index).
$ws.Cells.Item(45, 1).Value = 0

# Match the post-edit selection recorded in the workbook (activeCell /
# sqref both "B47").
$ws.Range("B47").Select()
